$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 284, pushing the existing
# rows 284-378 down to 286-380 (dimension grows from A1:R378 to A1:R380).
$ws.Range("A284:A285").EntireRow.Insert()

# Populate the new row 284 (Primera) with the new price-report entry.
$ws.Cells.Item(284, 1).Value = 1
$ws.Cells.Item(284, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(284, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(284, 4).Value = 44876
$ws.Cells.Item(284, 5).Value = 15
$ws.Cells.Item(284, 6).Value = 100112043
$ws.Cells.Item(284, 7).Value = "Pepino ensalada"
$ws.Cells.Item(284, 8).Value = "Sin especificar"
$ws.Cells.Item(284, 9).Value = "Primera"
$ws.Cells.Item(284, 10).Value = 150
$ws.Cells.Item(284, 11).Value = 15000
$ws.Cells.Item(284, 12).Value = 16000
$ws.Cells.Item(284, 13).Value = 15500
$ws.Cells.Item(284, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(284, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(284, 16).Value = 221
$ws.Cells.Item(284, 17).Value = 70
$ws.Cells.Item(284, 18).Value = "Hortaliza"

# Populate the new row 285 (Segunda) with the new price-report entry.
$ws.Cells.Item(285, 1).Value = 1
$ws.Cells.Item(285, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(285, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(285, 4).Value = 44876
$ws.Cells.Item(285, 5).Value = 15
$ws.Cells.Item(285, 6).Value = 100112043
$ws.Cells.Item(285, 7).Value = "Pepino ensalada"
$ws.Cells.Item(285, 8).Value = "Sin especificar"
$ws.Cells.Item(285, 9).Value = "Segunda"
$ws.Cells.Item(285, 10).Value = 160
$ws.Cells.Item(285, 11).Value = 11000
$ws.Cells.Item(285, 12).Value = 12000
$ws.Cells.Item(285, 13).Value = 11500
$ws.Cells.Item(285, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(285, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(285, 16).Value = 115
$ws.Cells.Item(285, 17).Value = 100
$ws.Cells.Item(285, 18).Value = "Hortaliza"
